$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.251.78"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.860.25"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'0.7019"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'237.58"
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.08258"
$ws.Range("E8").Value = "  +10.48%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").Value = "'0.08185"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.867.98"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.176"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.7137"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'89.14"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "29.267.28"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'0.000007855"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "'237.06"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "2.112.09"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'7.446"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'161.98"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'8.977"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'0.1444"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "'1.965"
$ws.Range("D30").Value = "'1.437"
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D31").Value = "'4.406"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "'4.065"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").Value = "'0.05210"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "'1.169"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("D36").Value = "'0.7075"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'1.003"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'2.668"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'0.01849"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").Value = "'0.9185"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "1.133.40"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("D43").Value = "'5.945"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'0.4282"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'70.69"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'102.41"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "2.009.95"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "'9.180"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  -0.83%  "
